$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.207771420478821
$ws.Range("B1").Value = 4.958170890808105
$ws.Range("C1").Value = 3.320107936859131
$ws.Range("D1").Value = 1.78874945640564
$ws.Range("E1").Value = 1.350434184074402
